# Översikt HANINGE.xlsx - automatic update of files
#
# This applies the following logical changes:
#  1. The two existing entries "A 67868-2021" (old row 5) and
#     "A 67894-2021" (old row 6) are removed from their old spot in
#     the table - they have been re-surveyed with updated
#     species-count / species-list data.
#  2. Refreshed versions of those same two entries are (re)inserted
#     at the very top of the data table, as the new rows 2 and 3.
#     Because of this, old rows 2-4 shift down to rows 4-6, while old
#     rows 7 onward keep their original row numbers (the removal of
#     old rows 5-6 and the insertion of the 2 new rows above cancel
#     out).
#  3. The "Förändrad" (column C) date stamp is bumped from 2023-09-19
#     (45188) to 2023-09-20 (45189) for every single data row in the
#     sheet, including the two new/refreshed rows above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Step 1: remove the old "A 67868-2021" / "A 67894-2021" rows (5
# and 6) from their current spot - their data is stale and is being
# replaced by the refreshed rows inserted at the top in step 2.
# ---------------------------------------------------------------
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()

# ---------------------------------------------------------------
# Step 2: insert two fresh blank rows above the current row 2 to
# hold the refreshed data for those same two entries.
# ---------------------------------------------------------------
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# ---------------------------------------------------------------
# Step 3: populate new row 2 -> "A 67868-2021" (refreshed)
# ---------------------------------------------------------------
$ws.Range("A2").Value = "A 67868-2021"
$ws.Range("B2").Value = 44525
$ws.Range("B2").NumberFormat = "YYYY-MM-DD"
$ws.Range("C2").Value = 45189
$ws.Range("C2").NumberFormat = "YYYY-MM-DD"
$ws.Range("D2").Value = "STOCKHOLMS LÄN"
$ws.Range("E2").Value = "HANINGE"
$ws.Range("G2").Value = 8.8
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 18
$ws.Range("R2").Value = "Gränsticka`r`nSpillkråka`r`nTallticka`r`nVedskivlav`r`nBjörksplintborre`r`nBlåmossa`r`nBronshjon`r`nFlagellkvastmossa`r`nGranbarkgnagare`r`nGrön sköldmossa`r`nKornknutmossa`r`nMindre märgborre`r`nPlatt fjädermossa`r`nScharlakansvårskål agg.`r`nStubbspretmossa`r`nSårläka`r`nBlåsippa`r`nRevlummer"
$ws.Range("R2").WrapText = $true
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/artfynd/A 67868-2021.xlsx", "A 67868-2021")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/kartor/A 67868-2021.png", "A 67868-2021")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/klagomål/A 67868-2021.docx", "A 67868-2021")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/klagomålsmail/A 67868-2021.docx", "A 67868-2021")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/tillsyn/A 67868-2021.docx", "A 67868-2021")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/tillsynsmail/A 67868-2021.docx", "A 67868-2021")'
$ws.Rows.Item(2).RowHeight = 15

# ---------------------------------------------------------------
# Step 4: populate new row 3 -> "A 67894-2021" (refreshed)
# ---------------------------------------------------------------
$ws.Range("A3").Value = "A 67894-2021"
$ws.Range("B3").Value = 44525
$ws.Range("B3").NumberFormat = "YYYY-MM-DD"
$ws.Range("C3").Value = 45189
$ws.Range("C3").NumberFormat = "YYYY-MM-DD"
$ws.Range("D3").Value = "STOCKHOLMS LÄN"
$ws.Range("E3").Value = "HANINGE"
$ws.Range("G3").Value = 5.1
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 9
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 18
$ws.Range("R3").Value = "Grönhjon`r`nKortskaftad ärgspik`r`nMotaggsvamp`r`nReliktbock`r`nSpillkråka`r`nSvartvit taggsvamp`r`nTallticka`r`nTalltita`r`nUllticka`r`nBjörksplintborre`r`nBlåmossa`r`nBronshjon`r`nGranbarkgnagare`r`nGrovticka`r`nGrön sköldmossa`r`nKornknutmossa`r`nStor revmossa`r`nStubbspretmossa"
$ws.Range("R3").WrapText = $true
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/artfynd/A 67894-2021.xlsx", "A 67894-2021")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/kartor/A 67894-2021.png", "A 67894-2021")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/klagomål/A 67894-2021.docx", "A 67894-2021")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/klagomålsmail/A 67894-2021.docx", "A 67894-2021")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/tillsyn/A 67894-2021.docx", "A 67894-2021")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/tillsynsmail/A 67894-2021.docx", "A 67894-2021")'
$ws.Rows.Item(3).RowHeight = 15

# ---------------------------------------------------------------
# Step 5: bump the "Förändrad" (column C) timestamp for every data
# row in the sheet from 2023-09-19 to 2023-09-20. Rows 2 and 3
# already carry the new date from the steps above, so re-stamping
# them is harmless.
# ---------------------------------------------------------------
$lastRow = $ws.UsedRange.Rows.Count()
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cur = $cell.Value()
    if ($cur -ne $null) {
        $cell.Value = 45189
        $cell.NumberFormat = "YYYY-MM-DD"
    }
}
